# Set up Instance Microservice Locally.docx
#
# The edit:
#   - The second paragraph (currently empty) gets a new run containing
#     "  *" (two spaces + an asterisk).
#   - The hidden "_GoBack" bookmark (Word's "last edit location" marker)
#     moves from its old, now-stale location (an empty paragraph further
#     down the document) to sit right after the newly typed text.
#
$d = $word.ActiveDocument

# 1) Remove the existing "_GoBack" bookmark from its old location, if present.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) Locate the (currently empty) second paragraph and replace its content
#    with a run of "  *" followed by a fresh "_GoBack" bookmark, while
#    preserving the paragraph's existing identity/formatting attributes.
$para = $d.Paragraphs.Item(2)
$range = $para.Range

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
         '<w:body>' +
           '<w:p w:rsidR="00DA6619" w:rsidRPr="00DA6619" w:rsidRDefault="00DA6619">' +
             '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
             '<w:r>' +
               '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
               '<w:t xml:space="preserve">  *</w:t>' +
             '</w:r>' +
             '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
             '<w:bookmarkEnd w:id="0"/>' +
           '</w:p>' +
         '</w:body>' +
       '</w:document>'

$range.InsertXML($xml)
